# TELEXo BOM update
#  - Row 15: fix the Device for the 180p cap (1206 package) from CAP0603-CAP -> CAP1206
#  - Merge the two "10k / Resistor 1%" rows (old rows 25 & 26) into a single
#    row: qty 4+8=12, keep the Mouser 71-CRCW0603-10K-E3 part, combine the
#    reference-designator lists from both rows (R1-R8 in F, R102-R105 in G).
#  - Remove the now-redundant second 10k row; everything below shifts up by one.
#  - Refresh the named range / selection to match the new 28-row extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: Device column fix ------------------------------------------------
$ws.Range("C15").Value = "CAP1206"

# --- Merge the duplicate 10k resistor rows into row 25 -----------------------
$ws.Range("A25").Value = 12
$ws.Range("F25").Value = "R1,R2,R3,R4,R5,R6,R7,R8"
$ws.Range("J25").Value = "71-CRCW0603-10K-E3"

# Drop the old row 26 (the second 10k line) - rows 27-29 shift up to 26-28,
# and the sheet dimension / table range update automatically.
$ws.Rows.Item(26).Delete()

# --- Keep the workbook-level defined name in sync with the new extent --------
$wb.Names.Item("Telex_O_1").RefersTo = "=TELEXo!`$A`$1:`$J`$28"

# --- Match the recorded selection state (whole row 25 selected) --------------
$ws.Rows.Item(25).Select() | Out-Null
